# This edit reorders the data rows of the worksheet (rows 2-25, columns A-F)
# so that each row's values move to a new row position, as if the
# underlying "totalBetInWins" data got added into the JSON payload and the
# rows were re-serialized in a different order. Row 1 (headers) and row 26
# (totals) are unchanged, as are rows 2, 5, 9, 11, 19, 24 and 25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of destination row -> new values for columns A..F
$rowValues = @{
    3  = @(1001, 18, 30, 75, 60, 72)
    4  = @(501, 9, 52, 30, 75, 45)
    6  = @(801, 3, 67, 65, 52, 45)
    7  = @(1202, 2, 10, 10, 10, 10)
    8  = @(901, 16, 15, 45, 60, 60)
    10 = @(1203, 3, 15, 15, 15, 15)
    12 = @(301, 6, 45, 30, 60, 45)
    13 = @(701, 3, 90, 45, 97, 15)
    14 = @(201, 9, 30, 15, 45, 30)
    15 = @(1201, 2, 10, 10, 10, 10)
    16 = @(802, 0, 4, 5, 4, 0)
    17 = @(2, 0, 2, 2, 2, 2)
    18 = @(502, 0, 4, 0, 0, 0)
    20 = @(3, 0, 3, 3, 3, 3)
    21 = @(1, 0, 2, 2, 2, 2)
    22 = @(402, 0, 0, 4, 0, 0)
    23 = @(602, 0, 0, 4, 0, 9)
}

$columns = @("A", "B", "C", "D", "E", "F")

foreach ($row in $rowValues.Keys) {
    $values = $rowValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$($row)").Value = $values[$i]
    }
}
